$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Date" column (B) timestamps for the existing MRF test case rows
# (B2:B7) to reflect the latest test run timestamps.
$ws.Range("B2").Value = "Fri Oct 25 12:05:10 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 12:05:23 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 12:05:35 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 12:05:47 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 12:05:59 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 12:06:11 EDT 2024"
